$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.736.03"
$ws.Range("E2").Value = "  +1.84%  "

# Row 3
$ws.Range("D3").Value = "2.208.69"
$ws.Range("E3").Value = "  -0.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'260.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.21%  "

# Row 6
$ws.Range("D6").Value = "'86.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.54%  "

# Row 7
$ws.Range("E7").Value = "  +0.74%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").Value = "'0.596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.51%  "

# Row 10
$ws.Range("D10").Value = "'45.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.36%  "

# Row 11
$ws.Range("D11").Value = "'0.0919"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.86%  "

# Row 12
$ws.Range("D12").Value = "'7.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.19%  "

# Row 13
$ws.Range("E13").Value = "  +1.60%  "

# Row 14
$ws.Range("D14").Value = "2.535.25"
$ws.Range("E14").Value = "  -0.35%  "

# Row 15
$ws.Range("D15").Value = "'14.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "

# Row 16
$ws.Range("D16").Value = "2.227.45"
$ws.Range("E16").Value = "  +0.45%  "

# Row 17
$ws.Range("D17").Value = "'0.784"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.84%  "

# Row 18
$ws.Range("D18").Value = "43.624.30"
$ws.Range("E18").Value = "  +1.86%  "

# Row 19
$ws.Range("E19").Value = "  +0.63%  "

# Row 20
$ws.Range("D20").Value = "'69.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.88%  "

# Row 21
$ws.Range("D21").Value = "'5.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "

# Row 22
$ws.Range("D22").Value = "'2.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.66%  "

# Row 23
$ws.Range("D23").Value = "'231.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.82%  "

# Row 24
$ws.Range("D24").Value = "'8.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.27%  "

# Row 26
$ws.Range("D26").Value = "'3.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.61%  "

# Row 27
$ws.Range("D27").Value = "'10.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.61%  "

# Row 28
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'39.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "

# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.92%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.46%  "

# Row 31
$ws.Range("D31").Value = "'174.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.88%  "

# Row 32
$ws.Range("D32").Value = "'20.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.97%  "

# Row 33
$ws.Range("D33").Value = "'0.0870"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.83%  "

# Row 34
$ws.Range("D34").Value = "'5.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.16%  "

# Row 35
$ws.Range("E35").Value = "  +1.40%  "

# Row 36
$ws.Range("E36").Value = "  +2.09%  "

# Row 37
$ws.Range("D37").Value = "'0.0358"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.20%  "

# Row 38
$ws.Range("E38").Value = "  +3.87%  "

# Row 39
$ws.Range("D39").Value = "'12.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.73%  "

# Row 40
$ws.Range("D40").Value = "'2.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.39%  "

# Row 41
$ws.Range("E41").Value = "  -0.24%  "

# Row 42
$ws.Range("D42").Value = "'63.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.72%  "

# Row 43
$ws.Range("D43").Value = "'5.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.14%  "

# Row 44
$ws.Range("D44").Value = "'0.199"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.36%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'100.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'8.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.33%  "

# Row 47
$ws.Range("D47").Value = "'0.0979"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "

# Row 48
$ws.Range("E48").Value = "  +4.32%  "

# Row 49
$ws.Range("E49").Value = "  +0.57%  "

# Row 50
$ws.Range("D50").Value = "'0.438"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.13%  "

# Row 51
$ws.Range("D51").Value = "'1.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.74%  "
